$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 379 (existing rows 379..473 shift down to 380..474)
$ws.Rows.Item(379).Insert()

# Populate the newly inserted row 379 with the new record
$ws.Cells.Item(379, 1).Value2 = 10
$ws.Cells.Item(379, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(379, 3).Value2 = "La Araucanía"
$ws.Cells.Item(379, 4).Value2 = 44736
$ws.Cells.Item(379, 5).Value2 = 9
$ws.Cells.Item(379, 6).Value2 = 100112032
$ws.Cells.Item(379, 7).Value2 = "Zapallo italiano"
$ws.Cells.Item(379, 8).Value2 = "Sin especificar"
$ws.Cells.Item(379, 9).Value2 = "Primera"
$ws.Cells.Item(379, 10).Value2 = 130
$ws.Cells.Item(379, 11).Value2 = 19000
$ws.Cells.Item(379, 12).Value2 = 20000
$ws.Cells.Item(379, 13).Value2 = 19385
$ws.Cells.Item(379, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(379, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(379, 16).Value2 = 323
$ws.Cells.Item(379, 17).Value2 = 60
$ws.Cells.Item(379, 18).Value2 = "Hortaliza"
